$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.671.63"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "2.967.06"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "594.29"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "145.41"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "2.966.91"
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("D9").Value = "0.505"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("E10").Value = "  +3.31%  "
$ws.Range("D11").Value = "0.146"
$ws.Range("E11").Value = "  +2.64%  "
$ws.Range("D12").Value = "0.445"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").Value = "0.0000238"
$ws.Range("E13").Value = "  +5.33%  "
$ws.Range("D14").Value = "33.13"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "3.460.11"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").Value = "62.605.78"
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("D18").Value = "6.71"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "2.971.29"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").Value = "441.57"
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("D21").Value = "13.47"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").Value = "0.671"
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("D23").Value = "7.07"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").Value = "81.68"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "11.23"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("E27").Value = "  -3.63%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  +3.68%  "
$ws.Range("D30").Value = "2.61"
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("E31").Value = "  -4.78%  "
$ws.Range("D32").Value = "0.0₃0953"
$ws.Range("E32").Value = "  +9.50%  "
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").Value = "26.47"
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").Value = "5.63"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").Value = "3.03"
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("D39").Value = "2.04"
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("D40").Value = "49.44"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("D41").Value = "8.54"
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("E42").Value = "  -4.72%  "
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").Value = "40.07"
$ws.Range("E44").Value = "  -4.60%  "
$ws.Range("D45").Value = "2.744.14"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").Value = "135.01"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").Value = "0.0340"
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("D48").Value = "361.67"
$ws.Range("E48").Value = "  -3.29%  "
$ws.Range("D50").Value = "23.03"
$ws.Range("E50").Value = "  -3.62%  "
$ws.Range("E51").Value = "  -0.48%  "
